# expansão das análises automáticas
# Adds three new columns (L: apoio_medio, M: contribuicoes, N: media_contribuicoes)
# with per-row data, matching the diff's expanded analysis table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns L, M, N
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# New data values per row (2-7)
$ws.Range("L2").Value = 92.57989061438856
$ws.Range("M2").Value = 229018
$ws.Range("N2").Value = 325.3096590909091

$ws.Range("L3").Value = 82.83823479360029
$ws.Range("M3").Value = 34535
$ws.Range("N3").Value = 274.0873015873016

$ws.Range("L4").Value = 89.86130176813208
$ws.Range("M4").Value = 177524
$ws.Range("N4").Value = 146.9569536423841

$ws.Range("L5").Value = 92.24386350483199
$ws.Range("M5").Value = 26122
$ws.Range("N5").Value = 149.2685714285714

$ws.Range("L6").Value = 19.47419260544111
$ws.Range("M6").Value = 2113
$ws.Range("N6").Value = 14.37414965986395

$ws.Range("L7").Value = 21.45251346829188
$ws.Range("M7").Value = 95
$ws.Range("N7").Value = 19
